$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.978.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.19%  '
$ws.Range("D3").Value = "'1.557.72"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.23%  '
$ws.Range("E4").Value = '  -0.28%  '
$ws.Range("D5").Value = "'206.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.13%  '
$ws.Range("E6").Value = '  +0.78%  '
$ws.Range("E7").Value = '  -0.30%  '
$ws.Range("D8").Value = "'22.12"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.87%  '
$ws.Range("E9").Value = '  -0.18%  '
$ws.Range("D10").Value = "'0.0596"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.07%  '
$ws.Range("D11").Value = "'0.0857"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.26%  '
$ws.Range("D12").Value = "'1.779.14"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.21%  '
$ws.Range("D13").Value = "'1.542.50"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.84%  '
$ws.Range("E14").Value = '  +0.95%  '
$ws.Range("E15").Value = '  +1.04%  '
$ws.Range("E16").Value = '  +0.24%  '
$ws.Range("D17").Value = "'26.976.05"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.26%  '
$ws.Range("E18").Value = '  +2.63%  '
$ws.Range("D19").Value = "'217.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.63%  '
$ws.Range("D20").Value = "'7.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.60%  '
$ws.Range("E21").Value = '  -0.30%  '
$ws.Range("E22").Value = '  +1.31%  '
$ws.Range("E23").Value = '  +0.59%  '
$ws.Range("E24").Value = '  -3.34%  '
$ws.Range("D25").Value = "'152.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.36%  '
$ws.Range("E26").Value = '  -0.24%  '
$ws.Range("E27").Value = '  +0.82%  '
$ws.Range("E28").Value = '  +1.24%  '
$ws.Range("E29").Value = '  -0.37%  '
$ws.Range("E30").Value = '  +1.18%  '
$ws.Range("E31").Value = '  +1.65%  '
$ws.Range("E32").Value = '  +0.36%  '
$ws.Range("D33").Value = "'1.421.38"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.40%  '
$ws.Range("D34").Value = "'3.11"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.77%  '
$ws.Range("D35").Value = "'1.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.52%  '
$ws.Range("D36").Value = "'1.04"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +7.89%  '
$ws.Range("E37").Value = '  +0.88%  '
$ws.Range("E38").Value = '  +0.33%  '
$ws.Range("D39").Value = "'0.531"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.57%  '
$ws.Range("E40").Value = '  +0.22%  '
$ws.Range("E41").Value = '  -0.31%  '
$ws.Range("D42").Value = "'5.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.05%  '
$ws.Range("E43").Value = '  +2.05%  '
$ws.Range("D44").Value = "'0.999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.22%  '
$ws.Range("D45").Value = "'64.73"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.71%  '
$ws.Range("D46").Value = "'1.75"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.29%  '
$ws.Range("D47").Value = "'1.692.62"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.17%  '
$ws.Range("D48").Value = "'87.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.50%  '
$ws.Range("E49").Value = '  +0.37%  '
$ws.Range("D50").Value = "'0.0₆0101"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.72%  '
$ws.Range("D51").Value = "'0.0959"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.09%  '
